$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.614.83'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '2.624.61'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.76'
$ws.Range("E5").Value = '  +2.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.91'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '2.630.72'
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E10").Value = '  -2.30%  '
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.334'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '3.080.83'
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").Value = '58.593.59'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.70'
$ws.Range("E16").Value = '  -2.22%  '
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").Value = '2.623.78'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '345.29'
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.43'
$ws.Range("E20").Value = '  -3.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.16'
$ws.Range("E21").Value = '  -2.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.12'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.37'
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.413'
$ws.Range("E25").Value = '  -2.03%  '
$ws.Range("E26").Value = '  +2.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("D28").Value = '0.0₃0799'
$ws.Range("E28").Value = '  -2.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.98'
$ws.Range("E29").Value = '  -1.01%  '
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.19'
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.83'
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.79'
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.974'
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.96'
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("E38").Value = '  +0.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.833'
$ws.Range("E39").Value = '  -3.66%  '
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.62'
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '278.76'
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0982'
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.52'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.595'
$ws.Range("E46").Value = '  -3.66%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0520'
$ws.Range("E47").Value = '  -3.58%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.30'
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("D50").Value = '1.977.90'
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("E51").Value = '  -2.60%  '
